# Commit message: "Cambiata gestione della colonna Anno in superdettagli"
#
# Semantic changes applied here:
#   1. Rename the (only) worksheet from "Sheet1" to "DoesNotMatter".
#      Excel automatically keeps the sheet-scoped defined name
#      (_xlnm._FilterDatabase) and any other sheet-qualified references
#      in sync with the new name.
#   2. Move the active-cell selection on that sheet from K12 to M27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet (this also updates the sheet-qualified
#    definedName "_xlnm._FilterDatabase" automatically).
$ws.Name = "DoesNotMatter"

# 2) Update the current selection / active cell on the sheet.
$ws.Range("M27").Select()
